# Update cryptocurrency Price (column D) and Volume(1h) (column E) figures
# to reflect the latest snapshot from the GitHub Actions refresh job.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Temporarily force column D to Text format so that purely-numeric-looking
# price strings (e.g. "315.03") are written back as text, matching the
# original inline-string cell type instead of being parsed as numbers.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "28.402.28"
$ws.Range("E2").Value = "  +0.13%  "
$ws.Range("D3").Value = "1.817.47"
$ws.Range("E3").Value = "  -0.47%  "
$ws.Range("E4").Value = "  +0.12%  "
$ws.Range("D5").Value = "315.03"
$ws.Range("E5").Value = "  -0.73%  "
$ws.Range("E6").Value = "  +0.07%  "
$ws.Range("D7").Value = "0.5130"
$ws.Range("E7").Value = "  -4.08%  "
$ws.Range("D8").Value = "0.3935"
$ws.Range("E8").Value = "  -3.01%  "
$ws.Range("D9").Value = "0.07964"
$ws.Range("E9").Value = "  +4.76%  "
$ws.Range("D10").Value = "41.67"
$ws.Range("E10").Value = "  -0.42%  "
$ws.Range("D11").Value = "1.107"
$ws.Range("E11").Value = "  +0.17%  "
$ws.Range("D12").Value = "20.98"
$ws.Range("E12").Value = "  +1.12%  "
$ws.Range("D13").Value = "6.248"
$ws.Range("E13").Value = "  -1.26%  "
$ws.Range("E14").Value = "  +0.08%  "
$ws.Range("D15").Value = "7.482"
$ws.Range("E15").Value = "  -1.67%  "
$ws.Range("D16").Value = "1.827.89"
$ws.Range("E16").Value = "  +0.19%  "
$ws.Range("D17").Value = "0.00001129"
$ws.Range("E17").Value = "  +5.22%  "
$ws.Range("D18").Value = "92.53"
$ws.Range("E18").Value = "  +3.60%  "
$ws.Range("D19").Value = "0.06625"
$ws.Range("E19").Value = "  +0.27%  "
$ws.Range("D20").Value = "17.66"
$ws.Range("E20").Value = "  +0.22%  "
$ws.Range("E21").Value = "  +0.02%  "
$ws.Range("D22").Value = "6.086"
$ws.Range("E22").Value = "  -0.33%  "
$ws.Range("D23").Value = "28.437.34"
$ws.Range("E23").Value = "  +0.21%  "
$ws.Range("D24").Value = "11.27"
$ws.Range("E24").Value = "  +0.71%  "
$ws.Range("D25").Value = "2.266"
$ws.Range("E25").Value = "  +3.75%  "
$ws.Range("D26").Value = "21.10"
$ws.Range("E26").Value = "  +2.64%  "
$ws.Range("D27").Value = "2.032.99"
$ws.Range("E27").Value = "  -0.13%  "
$ws.Range("D28").Value = "155.38"
$ws.Range("E28").Value = "  -1.56%  "
$ws.Range("D29").Value = "2.402"
$ws.Range("E29").Value = "  -2.69%  "
$ws.Range("D30").Value = "125.60"
$ws.Range("E30").Value = "  +1.35%  "
$ws.Range("D31").Value = "0.1101"
$ws.Range("E31").Value = "  +0.40%  "
$ws.Range("D32").Value = "1.101"
$ws.Range("E32").Value = "  -1.89%  "
$ws.Range("D33").Value = "5.697"
$ws.Range("E33").Value = "  +0.70%  "
$ws.Range("D34").Value = "3.652"
$ws.Range("E34").Value = "  +0.18%  "
$ws.Range("D35").Value = "0.07019"
$ws.Range("E35").Value = "  -3.82%  "
$ws.Range("D36").Value = "0.2219"
$ws.Range("E36").Value = "  -1.09%  "
$ws.Range("E37").Value = "  -0.97%  "
$ws.Range("D38").Value = "5.184"
$ws.Range("E38").Value = "  -0.10%  "
$ws.Range("D39").Value = "8.832"
$ws.Range("E39").Value = "  -0.56%  "
$ws.Range("D40").Value = "11.29"
$ws.Range("D41").Value = "0.6248"
$ws.Range("E41").Value = "  -0.09%  "
$ws.Range("D42").Value = "1.177"
$ws.Range("E42").Value = "  -0.75%  "
$ws.Range("D43").Value = "1.001"
$ws.Range("E43").Value = "  -0.05%  "
$ws.Range("E44").Value = "  +0.13%  "
$ws.Range("D45").Value = "13.49"
$ws.Range("D46").Value = "3.741"
$ws.Range("E46").Value = "  +0.95%  "
$ws.Range("D47").Value = "0.5899"
$ws.Range("E47").Value = "  +0.99%  "
$ws.Range("D48").Value = "124.85"
$ws.Range("E48").Value = "  -0.42%  "
$ws.Range("D49").Value = "1.972"
$ws.Range("E49").Value = "  -0.74%  "
$ws.Range("D50").Value = "1.187"
$ws.Range("E50").Value = "  -1.39%  "
$ws.Range("D51").Value = "0.06885"
$ws.Range("E51").Value = "  -0.05%  "

# Restore the original (default) cell style for column D so that no stray
# per-cell formatting is left behind.
$ws.Range("D2:D51").Style = "Normal"

